# Update the carjacking by-month YoY workbook for the new "through" date.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / tab to reflect the new cut-off date.
$ws.Name = "Through 2022-04-15"

# Update the label for the April row to reflect the new cut-off date.
$ws.Range("A5").Value = "April (through 04-15)"

# Update the April row's values (row 5).
$ws.Range("B5").Value = 11
$ws.Range("D5").Value = 28
$ws.Range("E5").Value = 26
$ws.Range("F5").Value = 26
$ws.Range("H5").Value = 52
$ws.Range("I5").Value = 62

# Update the Total row's values (row 6) accordingly.
$ws.Range("B6").Value = 77
$ws.Range("D6").Value = 217
$ws.Range("E6").Value = 223
$ws.Range("F6").Value = 136
$ws.Range("H6").Value = 475
$ws.Range("I6").Value = 496
